# Update status ("ESTADO") and last-modification date ("FECHA DE ULTIMA MODIFICACION")
# for the rows that were reprocessed on 04/01/2025.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = (Get-Date -Year 2025 -Month 1 -Day 4).Date

$rows = @(43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 77)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "PROCESADA"
    $ws.Range("C$r").Value = $newDate
}
